$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill data rows 2-44 with win/loss/tie totals
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 88   # AD
    $ws.Cells.Item($r, 31).Value = 74   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
